# Navigation for account tab created
# Adds a new "H2 console" error/solution entry to the "error report" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("error report")

# Widen column B so the longer text fits nicely.
$ws.Columns.Item(2).ColumnWidth = 105.45

# Fill in the Solution for the existing row 3 (the 'accessible' dependency error).
$ws.Range("C3").Value = "Opened H2 console before launching app"

# Add a brand new row (row 4) documenting the H2 connection error and its fix.
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 'org.h2.jdbc.JdbcSQLException: Connection is broken: "java.net.ConnectException: Connection refused: connect: localhost'
$ws.Range("C4").Value = "Restarted H2"

# Style B4's text like the error cells: red Arial 9pt.
# Build the font on a scratch named style (keeps the style table tidy - a
# single Font re-materialization instead of one per changed sub-property),
# apply it to the cell, then drop the temporary named style again.
$errStyle = $wb.Styles.Add("H2ErrorFont")
$errStyle.Font.Name = "Arial"
$errStyle.Font.Size = 9
$errStyle.Font.Color = 255
$ws.Range("B4").Style = "H2ErrorFont"
$wb.Styles.Item("H2ErrorFont").Delete()

# Move the active selection to the newly added cell, matching the edited view.
$ws.Range("C4").Select()
